# GitHub Actions symbol-list refresh for cryptos.xlsx (Thu Dec 22 15:24:22 UTC 2022)
# Updates the "Price" column for most rows, swaps the BitKan / HotbitToken rows
# (20 <-> 21) back to their earlier order, and refreshes a few "Volume(1h)" labels.
#
# The Price column is stored as TEXT in the workbook (e.g. "5.370" keeps its
# trailing zero), so every Price write below is prefixed with a leading
# apostrophe to force Excel to keep it as text instead of re-interpreting it
# as a number (which would drop trailing zeros / use scientific notation for
# the very small values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BNB
$ws.Range("D2").Value = "'241.91"
# OKB
$ws.Range("D3").Value = "'21.88"
# HuobiToken
$ws.Range("D4").Value = "'5.360"
# Cronos
$ws.Range("D5").Value = "'0.05709"
# GateToken
$ws.Range("D6").Value = "'3.423"
# MXToken
$ws.Range("D8").Value = "'0.8057"
# FTXToken
$ws.Range("D9").Value = "'0.8550"
# WazirX
$ws.Range("D10").Value = "'0.1439"
# MandalaExchangeToken
$ws.Range("D11").Value = "'0.07280"
# LiechtensteinCryptoassetsExchange
$ws.Range("D12").Value = "'0.03079"
# BitrueCoin
$ws.Range("D13").Value = "'0.03141"
# BitMartToken
$ws.Range("D14").Value = "'0.09355"
# MCDex
$ws.Range("D15").Value = "'3.936"
# BitForexToken
$ws.Range("D16").Value = "'0.001588"
# CoinExToken
$ws.Range("D17").Value = "'0.04824"
# One
$ws.Range("D18").Value = "'0.0005857"
# TigerCash
$ws.Range("D19").Value = "'0.006359"

# Row 20 becomes HotbitToken (was BitKan)
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004061"
$ws.Range("E20").Value = "19HotbitTokenHTB"

# Row 21 becomes BitKan (was HotbitToken)
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009956"
$ws.Range("E21").Value = "20BitKanKAN"

# NitroEx
$ws.Range("D22").Value = "'0.0001502"
# LEO
$ws.Range("D23").Value = "'3.723"
# BTSEToken
$ws.Range("D24").Value = "'2.175"
# BitpandaEcosystemToken
$ws.Range("D25").Value = "'0.3232"
# ProBitToken
$ws.Range("D26").Value = "'0.1242"
# UpBots
$ws.Range("D27").Value = "'0.0004005"

# IDEX
$ws.Range("D40").Value = "'0.03820"

# KickToken
$ws.Range("D41").Value = "'0.006776"
$ws.Range("E41").Value = "40KickTokenKICK"

# BKEXToken
$ws.Range("D42").Value = "'0.1049"
# CEJI
$ws.Range("D43").Value = "'0.003204"
# LocalTraders
$ws.Range("D44").Value = "'0.006615"
# Kangarootoken
$ws.Range("D46").Value = "'0.00000000751"

# CoinbaseStockToken
$ws.Range("D47").Value = "'0.5807"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

# BOLO
$ws.Range("D48").Value = "'0.1427"
# CryptobidCoin
$ws.Range("D49").Value = "'0.00002102"
# SpecialPowerGold
$ws.Range("D50").Value = "'0.01011"

Write-Host "Applied cryptos.xlsx updates"
